# Applies the "cryptos list" price/volume refresh described by the commit diff.
# For each changed row: Price (D) and Volume(1h) (E) are updated to the new scrape;
# rows 48/49 additionally swap their Coin/Link (B/C) because Bittensor/Monero
# traded ranking positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.240.20"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.708.62"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("D7").Value = "3.712.71"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.55"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "4.327.81"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "3.705.99"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "67.188.80"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.97"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.34%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "462.81"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.693"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000135"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -12.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "3.851.47"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.26"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.39"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.89"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").Value = "3.655.41"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -12.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.986"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.301"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "44.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "389.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.52%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0343"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.44%  "
